# Pendulum Team Hours Worked.xlsx
# Commit: "Added the time it took to do the status report"
#
# Matthew Handley's (row 4) hours for the "1/9 - 1/13" week (column D)
# increased from 1.25 to 1.75 to account for time spent on the status
# report. The "Total Hours/Week" row (row 9) contains a SUM formula over
# that column, so it recalculates automatically.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D4").Value = 1.75

# Leave the cursor where the user ended up after making the edit.
[void]$ws.Range("E13").Select()
